$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46061 -> 46062, i.e. 2026-02-08 -> 2026-02-09) for rows 2 through 24.
for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
